$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the documentation text in F6/F7 (Field::Pointer / Field::Type rows)
$ws.Range("F6").Value = "Field::Pointer gets set to nullptr_t, field.pointer gets set to nullptr (pointer type is used intermediately, but pointers to static references are not constexpr, hence setting nullptr)"
$ws.Range("F7").Value = "Field::Type gets set to decltype(&T::field)"

# Widen column F to fit the longer text, drop best-fit autosizing
$ws.Columns("F").ColumnWidth = 156

# Update the active selection
$ws.Range("F12").Select()
